# Final update to go client
# Re-sequence the product rows within several brand groups (Dinafex, Etorix,
# Flucloxin, Ketonic, Kynol, Zithrox): update Item Name (C), UOM (D) and
# TP (BB) so that each row reflects the re-ordered product list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dinafex group (rows 3-5): UOM stays the same ("30's") for all three rows.
$ws.Range("C3").Value = "Dinafex 180mg Tablet"
$ws.Range("BB3").Value = 224.89

$ws.Range("C4").Value = "Dinafex 120mg Tablet"
$ws.Range("BB4").Value = 179.91

$ws.Range("C5").Value = "Dinafex 60mg Tablet"
$ws.Range("BB5").Value = 78.70999999999999

# Etorix group (rows 7-9)
$ws.Range("C7").Value = "Etorix 120mg Tablet"
$ws.Range("D7").Value = "20's"

$ws.Range("C8").Value = "Etorix 60mg Tablet - 40's"
$ws.Range("D8").Value = "40's"
$ws.Range("BB8").Value = 209.9

$ws.Range("C9").Value = "Etorix 90mg Tablet"
$ws.Range("D9").Value = "30's"
$ws.Range("BB9").Value = 269.87

# Flucloxin group (rows 11-12)
$ws.Range("C11").Value = "Flucloxin 500mg Capsule"
$ws.Range("D11").Value = "30 's"
$ws.Range("BB11").Value = 237.74

$ws.Range("C12").Value = "Flucloxin 500mg Capsule - 36's"
$ws.Range("D12").Value = "36 's"
$ws.Range("BB12").Value = 284.21

# Ketonic group (rows 15-16)
$ws.Range("C15").Value = "Ketonic 10mg Tablet"
$ws.Range("D15").Value = "20's"
$ws.Range("BB15").Value = 150.38

$ws.Range("C16").Value = "Ketonic 30mg Injection"
$ws.Range("D16").Value = "5 's"
$ws.Range("BB16").Value = 206.77

# Kynol group (rows 17-19)
$ws.Range("C17").Value = "Kynol TR 200mg Capsule"
$ws.Range("D17").Value = "30 's"
$ws.Range("BB17").Value = 224.89

$ws.Range("C18").Value = "Kynol TR 100mg Capsule"
$ws.Range("D18").Value = "50 's"
$ws.Range("BB18").Value = 262.37

$ws.Range("C19").Value = "Kynol D 25mg Tablet"
$ws.Range("D19").Value = "60 's"
$ws.Range("BB19").Value = 180.45

# Row 20 (Naprox) keeps the same item name, only the UOM string changes
$ws.Range("D20").Value = "30 's"

# Zithrox group (rows 24-27)
$ws.Range("C24").Value = "Zithrox 15ml Suspension"
$ws.Range("D24").Value = "15 ml"
$ws.Range("BB24").Value = 71.95999999999999

$ws.Range("C25").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("D25").Value = "30ml"
$ws.Range("BB25").Value = 97.45

$ws.Range("C26").Value = "Zithrox 500mg Tablet"
$ws.Range("D26").Value = "6 's"
$ws.Range("BB26").Value = 136.83

$ws.Range("C27").Value = "Zithrox 250mg Tablet - 6's"
$ws.Range("D27").Value = "6's"
$ws.Range("BB27").Value = 89.95999999999999
